$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text before writing so that numeric-looking
# strings (e.g. "324.90", "0.5133") keep their literal text representation
# exactly as scraped, instead of being auto-coerced into floating point numbers.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

# Cell => new value updates, derived from the authoritative diff
$updates = @(
    @{ Cell = 'D2'; Value = '30.377.27' }
    @{ Cell = 'E2'; Value = '  +1.53%  ' }
    @{ Cell = 'D3'; Value = '2.011.04' }
    @{ Cell = 'E3'; Value = '  +4.39%  ' }
    @{ Cell = 'E4'; Value = '  +0.22%  ' }
    @{ Cell = 'D5'; Value = '324.90' }
    @{ Cell = 'E5'; Value = '  +1.47%  ' }
    @{ Cell = 'E6'; Value = '  +0.17%  ' }
    @{ Cell = 'D7'; Value = '0.5133' }
    @{ Cell = 'E7'; Value = '  +1.38%  ' }
    @{ Cell = 'D8'; Value = '0.4254' }
    @{ Cell = 'E8'; Value = '  +4.82%  ' }
    @{ Cell = 'D9'; Value = '0.08749' }
    @{ Cell = 'E9'; Value = '  +4.72%  ' }
    @{ Cell = 'B10'; Value = 'Polygon' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = 'D10'; Value = '1.134' }
    @{ Cell = 'E10'; Value = '  +2.58%  ' }
    @{ Cell = 'B11'; Value = 'Solana' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' }
    @{ Cell = 'D11'; Value = '24.45' }
    @{ Cell = 'E11'; Value = '  +2.57%  ' }
    @{ Cell = 'B12'; Value = 'WrappedEther' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D12'; Value = '2.015.47' }
    @{ Cell = 'E12'; Value = '  +5.02%  ' }
    @{ Cell = 'B13'; Value = 'Polkadot' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' }
    @{ Cell = 'D13'; Value = '6.632' }
    @{ Cell = 'E13'; Value = '  +3.35%  ' }
    @{ Cell = 'B14'; Value = 'Chainlink' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' }
    @{ Cell = 'D14'; Value = '7.452' }
    @{ Cell = 'E14'; Value = '  +2.76%  ' }
    @{ Cell = 'B15'; Value = 'BinanceUSD' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' }
    @{ Cell = 'D15'; Value = '1.005' }
    @{ Cell = 'E15'; Value = '  +0.48%  ' }
    @{ Cell = 'B16'; Value = 'Litecoin' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Cell = 'D16'; Value = '94.23' }
    @{ Cell = 'E16'; Value = '  +2.05%  ' }
    @{ Cell = 'B17'; Value = 'ShibaInu' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' }
    @{ Cell = 'D17'; Value = '0.00001113' }
    @{ Cell = 'E17'; Value = '  +1.25%  ' }
    @{ Cell = 'B18'; Value = 'TRON' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' }
    @{ Cell = 'D18'; Value = '0.06539' }
    @{ Cell = 'E18'; Value = '  +0.37%  ' }
    @{ Cell = 'B19'; Value = 'Avalanche' }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' }
    @{ Cell = 'D19'; Value = '18.84' }
    @{ Cell = 'E19'; Value = '  +2.99%  ' }
    @{ Cell = 'B20'; Value = 'Dai' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = 'D20'; Value = '1.001' }
    @{ Cell = 'E20'; Value = '  +0.15%  ' }
    @{ Cell = 'B21'; Value = 'Uniswap' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ Cell = 'D21'; Value = '6.206' }
    @{ Cell = 'E21'; Value = '  +4.17%  ' }
    @{ Cell = 'B22'; Value = 'WrappedBTC' }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' }
    @{ Cell = 'D22'; Value = '30.437.92' }
    @{ Cell = 'E22'; Value = '  +1.56%  ' }
    @{ Cell = 'B23'; Value = 'Cosmos' }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = 'D23'; Value = '11.84' }
    @{ Cell = 'E23'; Value = '  +4.44%  ' }
    @{ Cell = 'B24'; Value = 'Toncoin' }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Cell = 'D24'; Value = '2.254' }
    @{ Cell = 'E24'; Value = '  +2.94%  ' }
    @{ Cell = 'B25'; Value = 'WrappedliquidstakedEther2.0' }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Cell = 'D25'; Value = '2.254.18' }
    @{ Cell = 'E25'; Value = '  +5.25%  ' }
    @{ Cell = 'B26'; Value = 'EthereumClassic' }
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Cell = 'D26'; Value = '22.43' }
    @{ Cell = 'E26'; Value = '  +1.21%  ' }
    @{ Cell = 'B27'; Value = 'Monero' }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D27'; Value = '162.03' }
    @{ Cell = 'E27'; Value = '  -0.08%  ' }
    @{ Cell = 'B28'; Value = 'LidoDAOToken' }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = 'D28'; Value = '2.440' }
    @{ Cell = 'E28'; Value = '  +4.35%  ' }
    @{ Cell = 'B29'; Value = 'BitcoinCash' }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' }
    @{ Cell = 'D29'; Value = '131.34' }
    @{ Cell = 'E29'; Value = '  +1.68%  ' }
    @{ Cell = 'B30'; Value = 'ImmutableX' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D30'; Value = '1.141' }
    @{ Cell = 'E30'; Value = '  +0.71%  ' }
    @{ Cell = 'B31'; Value = 'Stellar' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' }
    @{ Cell = 'D31'; Value = '0.1055' }
    @{ Cell = 'E31'; Value = '  +1.69%  ' }
    @{ Cell = 'B32'; Value = 'Filecoin' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = 'D32'; Value = '6.089' }
    @{ Cell = 'E32'; Value = '  +1.80%  ' }
    @{ Cell = 'B33'; Value = 'HuobiToken' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ Cell = 'D33'; Value = '3.827' }
    @{ Cell = 'E33'; Value = '  +0.99%  ' }
    @{ Cell = 'B34'; Value = 'ARBITRUM' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = 'D34'; Value = '1.363' }
    @{ Cell = 'E34'; Value = '  +13.90%  ' }
    @{ Cell = 'B35'; Value = 'VeChain' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D35'; Value = '0.02537' }
    @{ Cell = 'E35'; Value = '  +3.50%  ' }
    @{ Cell = 'B36'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D36'; Value = '5.478' }
    @{ Cell = 'E36'; Value = '  +1.29%  ' }
    @{ Cell = 'B37'; Value = 'Hedera' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D37'; Value = '0.06673' }
    @{ Cell = 'E37'; Value = '  +3.64%  ' }
    @{ Cell = 'B38'; Value = 'Aptos' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = 'D38'; Value = '12.47' }
    @{ Cell = 'E38'; Value = '  +9.26%  ' }
    @{ Cell = 'B39'; Value = 'FraxShare' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D39'; Value = '9.195' }
    @{ Cell = 'E39'; Value = '  +5.00%  ' }
    @{ Cell = 'B40'; Value = 'Algorand' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Cell = 'D40'; Value = '0.2215' }
    @{ Cell = 'E40'; Value = '  +2.63%  ' }
    @{ Cell = 'B41'; Value = 'TheSandbox' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' }
    @{ Cell = 'D41'; Value = '0.6657' }
    @{ Cell = 'E41'; Value = '  +1.88%  ' }
    @{ Cell = 'B42'; Value = 'TrustWalletToken' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'D42'; Value = '1.238' }
    @{ Cell = 'E42'; Value = '  +1.61%  ' }
    @{ Cell = 'B43'; Value = 'Frax' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax' }
    @{ Cell = 'D43'; Value = '1.001' }
    @{ Cell = 'E43'; Value = '  +0.25%  ' }
    @{ Cell = 'B44'; Value = 'EnergySwap' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D44'; Value = '13.74' }
    @{ Cell = 'E44'; Value = '  +1.80%  ' }
    @{ Cell = 'B45'; Value = 'Decentraland' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D45'; Value = '0.6178' }
    @{ Cell = 'E45'; Value = '  +1.24%  ' }
    @{ Cell = 'B46'; Value = 'NEARProtocol' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Cell = 'D46'; Value = '2.199' }
    @{ Cell = 'E46'; Value = '  -1.80%  ' }
    @{ Cell = 'B47'; Value = 'PancakeSwap' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Cell = 'D47'; Value = '3.631' }
    @{ Cell = 'E47'; Value = '  -0.09%  ' }
    @{ Cell = 'B48'; Value = 'EOS' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos' }
    @{ Cell = 'D48'; Value = '1.259' }
    @{ Cell = 'E48'; Value = '  +3.97%  ' }
    @{ Cell = 'B49'; Value = 'Quant' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' }
    @{ Cell = 'D49'; Value = '124.93' }
    @{ Cell = 'E49'; Value = '  +2.29%  ' }
    @{ Cell = 'B50'; Value = 'Aave' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = 'D50'; Value = '81.21' }
    @{ Cell = 'E50'; Value = '  +2.68%  ' }
    @{ Cell = 'B51'; Value = 'Cronos' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = 'D51'; Value = '0.06915' }
    @{ Cell = 'E51'; Value = '  +1.56%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Drop the temporary Text number-format override so the cell style index
# reverts to the sheet default (matches the original, un-styled cells).
$textRange.ClearFormats()
